$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "filtrage recours par classe": add a classe_id foreign key field to the
# Coupon table so recourses (Recours) can be filtered by class.
# New field cell D10, matching the same field style used elsewhere (copy
# format from an existing field cell, e.g. C10, so it reuses the same
# shared fill style rather than creating a new one).
$ws.Range("C10").Copy()
$ws.Range("D10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D10").Value = "classe_id"

# Second new field cell, on a newly appended row (C12), same style.
$ws.Range("C10").Copy()
$ws.Range("C12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C12").Value = "classe_id"

$excel.CutCopyMode = $false

# Move the active selection to D11 (where the user left off editing).
$ws.Range("D11").Select()
